$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uni")

# Rotate the fuel-type codes in column A across the stacked ~TradeLinks blocks
$ws.Range("A2").Value = "FOL"
$ws.Range("A7").Value = "JET"
$ws.Range("A12").Value = "OTH"
$ws.Range("A17").Value = "COA"
$ws.Range("A22").Value = "PET"
$ws.Range("A27").Value = "COL"
$ws.Range("A32").Value = "DID"
$ws.Range("A37").Value = "DIJ"
$ws.Range("A42").Value = "LPG"
$ws.Range("A47").Value = "DSL"

# Move the trailing "1" flag from the COL/DID block's SI row (B34) to the
# JET/COA block's SI row (B19). Use an apostrophe-prefixed assignment so the
# value is stored as text (matching the other "1" flag cells, e.g. B29)
# rather than as a number, then reset the style so no quote-prefix
# formatting is left behind on the cell.
$ws.Range("B34").ClearContents()
$ws.Range("B19").Value = "'1"
$ws.Range("B19").Style = "Normal"
